# Append 8 new daily rows (2021-04-15 .. 2021-04-22) to the hospital
# ventilator dataset, continuing directly after the existing last row (415).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 416

$data = @(
    @("2021-04-15", 2196, 191, 1381, 624, 2780, 69, 489, 2222, 8.699999999999999, 62.89, 28.42, 2.48, 17.59, 79.93000000000001),
    @("2021-04-16", 2193, 193, 1361, 639, 2783, 70, 503, 2210, 8.800000000000001, 62.06, 29.14, 2.52, 18.07, 79.41),
    @("2021-04-17", 2188, 201, 1366, 621, 2781, 75, 483, 2223, 9.19, 62.43, 28.38, 2.7, 17.37, 79.94),
    @("2021-04-18", 2184, 200, 1266, 718, 2782, 76, 489, 2217, 9.16, 57.97, 32.88, 2.73, 17.58, 79.69),
    @("2021-04-19", 2175, 188, 1284, 703, 2781, 79, 475, 2227, 8.640000000000001, 59.03, 32.32, 2.84, 17.08, 80.08),
    @("2021-04-20", 2188, 197, 1346, 645, 2779, 76, 523, 2180, 9, 61.52, 29.48, 2.73, 18.82, 78.45),
    @("2021-04-21", 2197, 188, 1343, 666, 2779, 79, 524, 2176, 8.56, 61.13, 30.31, 2.84, 18.86, 78.3),
    @("2021-04-22", 2185, 185, 1333, 667, 2780, 78, 480, 2222, 8.470000000000001, 61.01, 30.53, 2.81, 17.27, 79.93000000000001)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    # Column A (DATE) holds a date-shaped string ("YYYY-MM-DD"). Assigning
    # it straight to .Value makes Excel auto-convert it to a date serial
    # number, but the source data stores it as plain text (a shared
    # string), same as every other row above it. Force the cell to Text
    # format before the assignment (the standard Excel technique for
    # entering digit/date-shaped literal text), then restore the default
    # "Normal" style afterwards so the cell ends up without any explicit
    # number format - matching the unstyled DATE cells elsewhere in the
    # column.
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $values[0]
    $dateCell.Style = "Normal"

    for ($c = 1; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
